# The data table gained one new observation which was inserted as row 207,
# pushing the existing rows 207-315 down to 208-316 (the former last row,
# 315, becomes row 316). We reproduce that by inserting a whole new row at
# position 207 (which shifts everything below it down automatically) and
# then populating the newly inserted row with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 207; rows 207..315 shift to 208..316.
$ws.Rows.Item(207).Insert()

# Populate the new row 207 with the new data record.
$ws.Range("A207").Value = 7
$ws.Range("B207").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C207").Value = "Ñuble"
$ws.Range("D207").Value = 45029
$ws.Range("E207").Value = 16
$ws.Range("F207").Value = 100112043
$ws.Range("G207").Value = "Pepino ensalada"
$ws.Range("H207").Value = "Sin especificar"
$ws.Range("I207").Value = "Primera"
$ws.Range("J207").Value = 120
$ws.Range("K207").Value = 11000
$ws.Range("L207").Value = 11000
$ws.Range("M207").Value = 11000
$ws.Range("N207").Value = '$/caja 80 unidades'
$ws.Range("O207").Value = "Región del Maule"
$ws.Range("P207").Value = 138
$ws.Range("Q207").Value = 80
$ws.Range("R207").Value = "Hortaliza"
